$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the existing last header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every data row (rows 2 through 46)
$ws.Range("AD2:AD46").Value = 76
$ws.Range("AE2:AE46").Value = 86
$ws.Range("AF2:AF46").Value = 0
